# Modificacion en consulta de productos y stock (codigo)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add "Agustina" / "en proceso" to the two rows that previously had no
# Responsable/Estado filled in.
$ws.Range("B48").Value = "Agustina"
$ws.Range("C48").Value = "en proceso"

$ws.Range("B49").Value = "Agustina"
$ws.Range("C49").Value = "en proceso"

# Move the active selection to C50 (single cell) as in the updated file.
$ws.Activate()
$ws.Range("C50").Select()
